$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 3: replace placeholder "None/None" prefix with real values
$ws.Range("A3").Value = "1002000/44444/Pan les"

# Row 4: add the missing label (I4/J4 already have values = 1)
$ws.Range("A4").Value = "1002000/44444/Pan yes"

# New row 5
$ws.Range("A5").Value = "1002000/44444/Pan hes"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("O5").Value = 30

# New row 6
$ws.Range("A6").Value = "1002000/44444/Pan kes"
$ws.Range("E6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("O6").Value = 30

# New row 7
$ws.Range("A7").Value = "1002000/44444/pan qes"
$ws.Range("G7").Value = 1

# New row 8
$ws.Range("A8").Value = "1002000/44444/pan bes"
$ws.Range("I8").Value = 1
